$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-11 Sunday", "2024-02-12 Monday"),
    @("471×9=4239", "703×3=2109"),
    @("182×4=728", "222×3=666"),
    @("676×8=5408", "800×2=1600"),
    @("396×5=1980", "144×7=1008"),
    @("401×5=2005", "152×9=1368"),
    @("169×2=338", "840×4=3360"),
    @("303×2=606", "430×2=860"),
    @("971×4=3884", "644×4=2576"),
    @("113×4=452", "183×3=549"),
    @("665×8=5320", "636×4=2544"),
    @("885×8=7080", "239×5=1195"),
    @("600×2=1200", "295×5=1475"),
    @("441×3=1323", "677×9=6093"),
    @("877×7=6139", "821×6=4926"),
    @("179×4=716", "759×9=6831"),
    @("145×4=580", "365×8=2920"),
    @("845×8=6760", "772×8=6176"),
    @("944×5=4720", "602×6=3612"),
    @("743×8=5944", "742×2=1484"),
    @("180×7=1260", "270×3=810"),
    @("565×5=2825", "891×8=7128"),
    @("215×9=1935", "886×5=4430"),
    @("461×6=2766", "156×2=312"),
    @("558×6=3348", "522×9=4698"),
    @("343×6=2058", "149×4=596")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
